# Fixed a bug in IVal(str)
# The data rows (2-25) on the active sheet got reshuffled because the
# sort/parse routine that turns the "symbol" text into a numeric key
# (IVal) was buggy. This re-applies the corrected row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(401, 9, 48, 67, 75, 45),
    @(201, 9, 30, 15, 45, 30),
    @(901, 16, 15, 45, 60, 60),
    @(1001, 18, 30, 75, 60, 72),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(601, 9, 60, 67, 60, 42),
    @(801, 3, 67, 65, 52, 45),
    @(101, 9, 30, 15, 60, 15),
    @(501, 9, 52, 30, 75, 45),
    @(902, 1, 0, 0, 0, 0),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(802, 0, 4, 5, 4, 0),
    @(1, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(1101, 0, 15, 30, 30, 0),
    @(402, 0, 0, 4, 0, 0),
    @(602, 0, 0, 4, 0, 9),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
